# Implemented decode and encode to base64
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Header row: keep the same labels ("usuario" / "contraseña"), but make
# them bold and horizontally centered.
$headerCell = $ws.Cells.Item(1, 1)
$headerCell.Value = "usuario"
$headerCell.HorizontalAlignment = -4108  # xlCenter
$headerCell.Font.Bold = $true

$headerCell.Copy()
$secondHeaderCell = $ws.Cells.Item(1, 2)
$secondHeaderCell.Value = "contraseña"
$secondHeaderCell.PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false

# Replace the placeholder "prueba" rows with base64-encoded user/password
# pairs (decode/encode to base64).
$ws.Cells.Item(2, 1).Value = "dXN1YXJpb1BydWViYQ=="
$ws.Cells.Item(2, 2).Value = "MWEyUzNkNEZ0ZXN0UGFzc3dvcmQ="

$ws.Cells.Item(3, 1).Value = "dXN1YXJpb1Bpenph"
$ws.Cells.Item(3, 2).Value = "MXEyVzNlNFJ0ZXN0UGFzc3dvcmQ="

$ws.Cells.Item(4, 1).Value = "dXN1YXJpb0NhcGliYXJh"
$ws.Cells.Item(4, 2).Value = "OU04bjdCNnZ0ZXN0UGFzc3dvcmQ="
